$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 22, shifting all subsequent rows (old 22..133) down to (23..134).
$ws.Rows.Item(22).Insert()

# Populate the newly inserted row 22 with the new weekly record.
$ws.Range("A22").Value = 11
$ws.Range("B22").Value = "Vega Monumental Concepción"
$ws.Range("C22").Value = "Bíobío"
$ws.Range("D22").Value = (Get-Date -Year 2022 -Month 1 -Day 14 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E22").Value = 8
$ws.Range("F22").Value = "Fruta"
$ws.Range("G22").Value = 100108
$ws.Range("H22").Value = "Tropicales y subtropicales"
$ws.Range("I22").Value = 100108005
$ws.Range("J22").Value = "Piña"
$ws.Range("K22").Value = "Caramelo"
$ws.Range("L22").Value = "Segunda"
$ws.Range("M22").Value = 200
$ws.Range("N22").Value = 16000
$ws.Range("O22").Value = 16500
$ws.Range("P22").Value = 16250
$ws.Range("Q22").Value = "$/caja 14 unidades"
$ws.Range("R22").Value = "Ecuador"
$ws.Range("S22").Value = 1161
$ws.Range("T22").Value = 14
